$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:G across rows 2-11 (ifoCAST full series evaluation)
$data = @(
    @(-0.3731803712827915, 2.007670690660352, 12.87239756141505, 3.58781236429876, 3.652324422179951, 22),
    @(-0.4804827079566977, 1.161321962968468, 5.130197931044437, 2.264994024505238, 2.268105159228511, 21),
    @(-0.2697575283137488, 0.8727440841527818, 2.711211561858076, 1.64657570790355, 1.666525676141245, 20),
    @(-0.02251358845263601, 0.8231985852708921, 1.758851113115051, 1.326216842418709, 1.362361935462548, 19),
    @(-0.1079376821351743, 0.7903692723160057, 1.18367369466905, 1.087967690085073, 1.113986395682008, 18),
    @(0.05847297744281363, 0.6864530745518149, 0.8768293907560101, 0.9363916866119701, 0.9633267594741898, 17),
    @(0.1260832718587093, 0.5907394480024477, 0.5291801755352681, 0.7274477132655432, 0.7399338234140823, 16),
    @(0.189738633461222, 0.6135846084182608, 0.6359933741545197, 0.7974919273287472, 0.8017787981776971, 15),
    @(0.2282491029459698, 0.5916736138517826, 0.4825523572072717, 0.6946598859926142, 0.6808572570925593, 14),
    @(0.2349671465221997, 0.5153392861830149, 0.3670491565632086, 0.6058458191348758, 0.5812282079097323, 13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
